# Apply updated yearly financial figures to the REDU sheet (Income
# Statement / Balance Sheet / Cash Flow Statement tables).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REDU")

# Row 8: Total Revenue
$ws.Range("D8").Value = 143900
$ws.Range("E8").Value = 105500
$ws.Range("F8").Value = 78600
$ws.Range("G8").Value = 60400

# Row 9: Cost of Revenue
$ws.Range("D9").Value = 67100
$ws.Range("E9").Value = 54000
$ws.Range("F9").Value = 51400
$ws.Range("G9").Value = 43800

# Row 10: Gross Profit
$ws.Range("D10").Value = 76700
$ws.Range("E10").Value = 51600
$ws.Range("F10").Value = 27100
$ws.Range("G10").Value = 16600

# Row 17: Total Operating Expenses
$ws.Range("D17").Value = 143900
$ws.Range("E17").Value = 95000
$ws.Range("F17").Value = 85900
$ws.Range("G17").Value = 73100

# Row 18: Operating Income or Loss
$ws.Range("E18").Value = 10500
$ws.Range("F18").Value = -7300
$ws.Range("G18").Value = -12700

# Row 20: Total Other Income/Expenses Net
$ws.Range("F20").Value = 2500
$ws.Range("G20").Value = 1100

# Row 21: Earnings Before Interest And Taxes
$ws.Range("D21").Value = 14200
$ws.Range("E21").Value = 26800
$ws.Range("F21").Value = 10000
$ws.Range("G21").Value = "NA"

# Row 23: Income Before Tax
$ws.Range("E23").Value = 12300
$ws.Range("F23").Value = -4900
$ws.Range("G23").Value = -11600

# Row 24: Income Tax Expense
$ws.Range("D24").Value = 7900
$ws.Range("E24").Value = 4800

# Row 26: Income After Tax
$ws.Range("D26").Value = -8000
$ws.Range("E26").Value = 7500
$ws.Range("F26").Value = -4700
$ws.Range("G26").Value = -10800

# Row 27: Net Income From Continuing Ops
$ws.Range("D27").Value = -7100
$ws.Range("E27").Value = 8000
$ws.Range("F27").Value = -3900
$ws.Range("G27").Value = -9700

# Row 32: Other Items
$ws.Range("F32").Value = -2500
$ws.Range("G32").Value = -1100

# Row 33: Net Income
$ws.Range("D33").Value = -7100
$ws.Range("E33").Value = 8000
$ws.Range("F33").Value = -3900
$ws.Range("G33").Value = -9700

# Row 35: Net Income Applicable To Common Shares
$ws.Range("D35").Value = -7100
$ws.Range("E35").Value = 8000
$ws.Range("F35").Value = -3900
$ws.Range("G35").Value = -9700

# Row 41: Cash And Cash Equivalents
$ws.Range("D41").Value = 156700
$ws.Range("E41").Value = 95000
$ws.Range("F41").Value = 76800

# Row 43: Net Receivables
$ws.Range("D43").Value = 1700

# Row 44: Inventory
$ws.Range("D44").Value = 1200

# Row 45: Other Current Assets
$ws.Range("D45").Value = 10000
$ws.Range("E45").Value = 9100
$ws.Range("F45").Value = 4200

# Row 46: Total Current Assets
$ws.Range("D46").Value = 169600
$ws.Range("E46").Value = 105000
$ws.Range("F46").Value = 82100

# Row 48: Property Plant and Equipment
$ws.Range("D48").Value = 14900
$ws.Range("E48").Value = 11200
$ws.Range("F48").Value = 10500

# Row 49: Goodwill
$ws.Range("D49").Value = 100400
$ws.Range("E49").Value = 102100
$ws.Range("F49").Value = 102300

# Row 52: Other Assets
$ws.Range("D52").Value = 5500

# Row 54: Total Assets
$ws.Range("D54").Value = 290300
$ws.Range("E54").Value = 222700
$ws.Range("F54").Value = 198200

# Row 57: Accounts Payable
$ws.Range("D57").Value = 4700

# Row 58: Short/Current Long Term Debt
$ws.Range("E58").Value = 5700

# Row 59: Other Current Liabilities
$ws.Range("D59").Value = 148300
$ws.Range("E59").Value = 107000
$ws.Range("F59").Value = 84400

# Row 60: Total Current Liabilities
$ws.Range("D60").Value = 153000
$ws.Range("E60").Value = 113300
$ws.Range("F60").Value = 84800

# Row 61: Long Term Debt
$ws.Range("D61").Value = 92500
$ws.Range("E61").Value = 49400

# Row 62: Other Liabilities
$ws.Range("D62").Value = 1000

# Row 66: Total Liabilities
$ws.Range("D66").Value = 244300
$ws.Range("E66").Value = 162200
$ws.Range("F66").Value = 85900

# Row 72: Retained Earnings
$ws.Range("D72").Value = -39900
$ws.Range("E72").Value = -15100
$ws.Range("F72").Value = -23100

# Row 76: Total Stockholder Equity
$ws.Range("D76").Value = 46000
$ws.Range("E76").Value = 60400
$ws.Range("F76").Value = 112300

# Row 81: Net Income
$ws.Range("D81").Value = -7100
$ws.Range("E81").Value = 8000
$ws.Range("F81").Value = -3900
$ws.Range("G81").Value = -9700

# Row 83: Depreciation
$ws.Range("D83").Value = 10400
$ws.Range("E83").Value = 13600
$ws.Range("F83").Value = 14900
$ws.Range("G83").Value = "NA"

# Row 89: Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 35600
$ws.Range("E89").Value = 24300
$ws.Range("F89").Value = 3700
$ws.Range("G89").Value = "NA"

# Row 91: Capital Expenditures
$ws.Range("D91").Value = -5300
$ws.Range("E91").Value = -5200
$ws.Range("F91").Value = -5400
$ws.Range("G91").Value = "NA"

# Row 94: Total Cash Flows From Investing Activities
$ws.Range("D94").Value = -6300
$ws.Range("E94").Value = -5700
$ws.Range("F94").Value = -6900
$ws.Range("G94").Value = "NA"

# Row 100: Total Cash Flows From Financing Activities
$ws.Range("D100").Value = -11900
$ws.Range("E100").Value = 0
$ws.Range("F100").Value = 600
$ws.Range("G100").Value = "NA"

# Row 101: Effect Of Exchange Rate Changes 
$ws.Range("D101").Value = 800
$ws.Range("E101").Value = 400
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = "NA"

# Row 102: Change In Cash and Cash Equivalents 
$ws.Range("D102").Value = 18200
$ws.Range("E102").Value = 19000
$ws.Range("F102").Value = -2600
$ws.Range("G102").Value = "NA"
